# Applies the "finished maze gen/added cam movement" update to the DTT hour log.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing log entry description (row 6 / D6): append text about
#     backtracking & wall removal now that the algorithm is finished.
$ws.Range("D6").Value = "it now loops trough the nodes and selects one randomly that it will chose as the next node and it now backtracks and removes walls (algorithm done)"

# --- Hours logged on the existing two entries were corrected from 2 to 1.
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1

# --- New log entry in row 28: "Sources in the research document".
$ws.Range("A28").Value = "Sources in the research document"

# --- New log entry in row 7: "implementing userstory 2" on 27 Nov 2022.
$ws.Range("A7").Value = "implementing userstory 2"
$ws.Range("C7").Value = "11/27/2022"

# --- Move the active selection to F7 (matches the author's last cursor spot).
$ws.Range("F7").Select()

$wb.Save()
